$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values (post-edit) for columns D, J, K, L, M, N, O, P, Q for rows 2-13.
# The edit re-distributes the existing per-row records across rows (a weekly
# re-sequencing of the daily records), leaving columns A, B, C, E, F, G, H, I, R unchanged.

$data = @{
    2  = @{ D = 44315; J = 25; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 167; Q = 60 }
    3  = @{ D = 44312; J = 30; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 167; Q = 60 }
    4  = @{ D = 44186; J = 15; K = 7000;  L = 7000;  M = 7000;  N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 117; Q = 60 }
    5  = @{ D = 44284; J = 35; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 167; Q = 60 }
    6  = @{ D = 44277; J = 25; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 167; Q = 60 }
    7  = @{ D = 44333; J = 25; K = 10000; L = 11000; M = 10400; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 173; Q = 60 }
    8  = @{ D = 44405; J = 45; K = 9000;  L = 9000;  M = 9000;  N = "`$/caja 50 unidades"; O = "Provincia de Quillota"; P = 180; Q = 50 }
    9  = @{ D = 44243; J = 80; K = 10000; L = 11000; M = 10375; N = "`$/caja 60 unidades"; O = "Provincia de Quillota"; P = 173; Q = 60 }
    10 = @{ D = 44585; J = 30; K = 11000; L = 11000; M = 11000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 183; Q = 60 }
    11 = @{ D = 44200; J = 10; K = 9000;  L = 9000;  M = 9000;  N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 150; Q = 60 }
    12 = @{ D = 44179; J = 15; K = 7000;  L = 7000;  M = 7000;  N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 117; Q = 60 }
    13 = @{ D = 44291; J = 20; K = 9000;  L = 9000;  M = 9000;  N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 150; Q = 60 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D    # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals.L   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 14).Value = $vals.N   # N - Unidad de comercializacion
    $ws.Cells.Item($row, 15).Value = $vals.O   # O - Origen
    $ws.Cells.Item($row, 16).Value = $vals.P   # P - Precio $/Kg
    $ws.Cells.Item($row, 17).Value = $vals.Q   # Q - Kg o Unidades
}
